$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete entire row 794 (the "眠たい目つきと暖かい微笑みで" post), shifting subsequent rows up.
$ws.Rows.Item(794).Delete()
